# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new labels in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing bold/centered/bordered header style (style index 1)
# by copying format from an existing header cell instead of re-deriving it.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Data rows: every player row (2-46) gets the same season record values.
for ($r = 2; $r -le 46; $r++) {
    $ws.Range("AD$r").Value = 72
    $ws.Range("AE$r").Value = 90
    $ws.Range("AF$r").Value = 0
}

$excel.CutCopyMode = 0
